# Auto-generated edit script: update cryptos list values (Thu Jul  4 04:58:32 UTC 2024)
$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

function Set-TextValue($cell, $text) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

Set-TextValue "D2" "59.015.94"
Set-TextValue "E2" "  -3.00%  "
Set-TextValue "D3" "3.237.14"
Set-TextValue "E3" "  -3.61%  "
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "541.03"
Set-TextValue "E5" "  -4.57%  "
Set-TextValue "D6" "136.87"
Set-TextValue "E6" "  -7.63%  "
Set-TextValue "E7" "  -0.18%  "
Set-TextValue "D8" "3.235.49"
Set-TextValue "E8" "  -3.66%  "
Set-TextValue "E9" "  -4.10%  "
Set-TextValue "D10" "7.65"
Set-TextValue "E10" "  -3.54%  "
Set-TextValue "E11" "  -5.26%  "
Set-TextValue "E12" "  -4.02%  "
Set-TextValue "D13" "3.788.80"
Set-TextValue "E13" "  -3.65%  "
Set-TextValue "E14" "  -1.03%  "
Set-TextValue "D15" "26.04"
Set-TextValue "E15" "  -6.83%  "
Set-TextValue "D16" "3.231.39"
Set-TextValue "E16" "  -3.61%  "
Set-TextValue "E17" "  -5.43%  "
Set-TextValue "D18" "59.033.60"
Set-TextValue "E18" "  -3.18%  "
Set-TextValue "D19" "5.91"
Set-TextValue "E19" "  -6.83%  "
Set-TextValue "D20" "13.35"
Set-TextValue "E20" "  -5.69%  "
Set-TextValue "E21" "  -6.07%  "
Set-TextValue "D22" "362.88"
Set-TextValue "E22" "  -2.87%  "
Set-TextValue "E23" "  +0.05%  "
Set-TextValue "D24" "70.59"
Set-TextValue "E24" "  -6.33%  "
Set-TextValue "D25" "0.522"
Set-TextValue "E25" "  -6.65%  "
Set-TextValue "D26" "3.367.93"
Set-TextValue "E26" "  -3.81%  "
Set-TextValue "B27" "Kaspa"
Set-TextValue "C27" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D27" "0.171"
Set-TextValue "E27" "  -2.68%  "
Set-TextValue "B28" "PEPE"
Set-TextValue "C28" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D28" "0.0₃0976"
Set-TextValue "E28" "  -9.97%  "
Set-TextValue "E29" "  +0.13%  "
Set-TextValue "D30" "7.11"
Set-TextValue "E30" "  -3.48%  "
Set-TextValue "E31" "  +0.01%  "
Set-TextValue "E32" "  -6.47%  "
Set-TextValue "D33" "7.13"
Set-TextValue "E33" "  -7.10%  "
Set-TextValue "D34" "21.96"
Set-TextValue "E34" "  -3.82%  "
Set-TextValue "E35" "  -4.83%  "
Set-TextValue "E36" "  -7.42%  "
Set-TextValue "D37" "163.08"
Set-TextValue "E37" "  -3.21%  "
Set-TextValue "D38" "6.44"
Set-TextValue "E38" "  -4.90%  "
Set-TextValue "E39" "  -6.43%  "
Set-TextValue "D40" "26.49"
Set-TextValue "E40" "  -9.33%  "
Set-TextValue "E41" "  -5.04%  "
Set-TextValue "D42" "3.268.11"
Set-TextValue "E42" "  -3.68%  "
Set-TextValue "D43" "41.25"
Set-TextValue "E43" "  -2.45%  "
Set-TextValue "D44" "0.719"
Set-TextValue "E44" "  -5.33%  "
Set-TextValue "D45" "1.12"
Set-TextValue "E45" "  -2.23%  "
Set-TextValue "E46" "  -5.63%  "
Set-TextValue "E48" "  -0.02%  "
Set-TextValue "D49" "2.306.40"
Set-TextValue "E49" "  -7.30%  "
Set-TextValue "D50" "6.32"
Set-TextValue "E50" "  -5.31%  "
Set-TextValue "D51" "21.03"
Set-TextValue "E51" "  -6.64%  "
